$d = $word.ActiveDocument

# Locate the paragraph that ends with "Large Class" (the last existing
# code-smell entry in section 1.2) so we can insert the new code-smell
# entries right after it, before the pre-existing blank spacer paragraphs.
$findRange = $d.Content
$found = $findRange.Find.Execute("Large Class", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Large Class' paragraph"
}

$largeClassPara = $findRange.Paragraphs(1)
$anchorPara = $largeClassPara.Next()

$insertionXml = @'
<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="3"/>
        </w:numPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t xml:space="preserve">Where: </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="SubtitleChar"/>
        </w:rPr>
        <w:t>PrinterManager</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> class, the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="SubtitleChar"/>
        </w:rPr>
        <w:t>addPrint</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="SubtitleChar"/>
        </w:rPr>
        <w:t>()</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> method</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="3"/>
        </w:numPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t xml:space="preserve">What’s wrong: </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">The method accepts 7 different parameters. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="3"/>
        </w:numPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>Why is it wrong:</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Similar to previous “code smells”, this large </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>amount</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> of parameters makes it harder to read and, of course, maintain.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="3"/>
        </w:numPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t xml:space="preserve">Possible solution: </w:t>
      </w:r>
      <w:r>
        <w:t>Accept one parameter of Print class</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="3"/>
        </w:numPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t xml:space="preserve">Name of code smell: </w:t>
      </w:r>
      <w:r>
        <w:t>Long Parameter List</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="3"/>
        </w:numPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t xml:space="preserve">Where: </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="SubtitleChar"/>
        </w:rPr>
        <w:t>PrinterManager</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> class, the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="SubtitleChar"/>
        </w:rPr>
        <w:t>getSpoolByID</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="SubtitleChar"/>
        </w:rPr>
        <w:t>()</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> method</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="3"/>
        </w:numPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t xml:space="preserve">What’s wrong: </w:t>
      </w:r>
      <w:r>
        <w:t>The m</w:t>
      </w:r>
      <w:r>
        <w:t>ethod is not used in the program</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="3"/>
        </w:numPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>Why is it wrong:</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>It is unnecessary to keep unused methods</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="3"/>
        </w:numPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t xml:space="preserve">Possible solution: </w:t>
      </w:r>
      <w:r>
        <w:t>Remove method</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="3"/>
        </w:numPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t xml:space="preserve">Name of code smell: </w:t>
      </w:r>
      <w:r>
        <w:t>Dead code</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="3"/>
        </w:numPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="3"/>
        </w:numPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t xml:space="preserve">Where: </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="SubtitleChar"/>
        </w:rPr>
        <w:t>PrinterManager</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> class, the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="SubtitleChar"/>
        </w:rPr>
        <w:t>registerCompletion</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="SubtitleChar"/>
        </w:rPr>
        <w:t>()</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">and </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="SubtitleChar"/>
        </w:rPr>
        <w:t>register</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="SubtitleChar"/>
        </w:rPr>
        <w:t>PrintFailure</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rStyle w:val="SubtitleChar"/>
        </w:rPr>
        <w:t>()</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>method</w:t>
      </w:r>
      <w:r>
        <w:t>s</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="3"/>
        </w:numPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t xml:space="preserve">What’s wrong: </w:t>
      </w:r>
      <w:r>
        <w:t>Both methods contain similar code for finding a printer by ID and handling the print task.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="3"/>
        </w:numPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t>Why is it wrong:</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>D</w:t>
      </w:r>
      <w:r>
        <w:t>uplicate code increases maintenance effort and the risk of inconsistencies</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="3"/>
        </w:numPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t xml:space="preserve">Possible solution: </w:t>
      </w:r>
      <w:r>
        <w:t>Extract duplicate code into a separate method</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="3"/>
        </w:numPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t xml:space="preserve">Name of code smell: </w:t>
      </w:r>
      <w:r>
        <w:t>Duplicate Code</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
</pkg:xmlData>
'@

$target = $d.Range($anchorPara.Range.Start, $anchorPara.Range.Start)
$result = $target.InsertXML($insertionXml)

# InsertXML drops <w:rStyle> (character style) references from the inline
# <w:rPr> it inserts, so re-apply the "Subtitle Char" character style to the
# method / class-name mentions that need it (matches the source XML intent).
$styleTargets = @("PrinterManager", "addPrint()", "PrinterManager", "getSpoolByID()", "PrinterManager", "registerCompletion()", "registerPrintFailure()")
$searchStart = $anchorPara.Range.Start
foreach ($term in $styleTargets) {
    $fr = $d.Range($searchStart, $d.Content.End)
    $ok = $fr.Find.Execute($term, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find style target: $term"
    }
    $fr.Style = "Subtitle Char"
    $searchStart = $fr.End
}

Write-Output "Inserted new code-smell paragraphs."
Write-Output $d.Paragraphs.Count
